# Automatische test-sync: 2025-06-26 19:37:50
# Adds the new "Logs" row (row 12) captured from the MailMind test mailbox,
# extends the conditional-formatting ranges to cover it, and bumps the
# "Bestelling / Levering" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 12 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Kun je 2 dozen nitrilhandschoenen bestellen?"
$logs.Range("B12").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C12").Value = "He Johan,`nKun je 2 dozen nitrilhandschoenen bestellen?`nMarc`nSent using {0}"
$logs.Range("D12").Value = "Bestelling / Levering"
$logs.Range("E12").Value = "Bedankt voor je bericht. Ik neem dit z.s.m. in behandeling."
$logs.Range("F12").Value = "2025-06-26 19:37:15"
$logs.Range("G12").Value = "Ja"
$logs.Range("H12").Value = "Ja"
$logs.Range("I12").Value = "Nee"

# --- Extend conditional formatting ranges to include the new row ---------
$logs.Range("D2:D11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D12"))
$logs.Range("G2:G11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G12"))
$logs.Range("H2:H11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H12"))
$logs.Range("I2:I11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I12"))

# --- Dashboard sheet: bump the "Bestelling / Levering" count -------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 7
